# Configuration file update: use "SGL" instead of "Dbl" for floating point
# data-type values in the "Sample csv" sheet, and leave the selection on
# the cell the author was last working on (C35).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$used = $ws.UsedRange
$rowCount = $used.Rows.Count
$colCount = $used.Columns.Count

for ($r = 1; $r -le $rowCount; $r++) {
    for ($c = 1; $c -le $colCount; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        if ($cell.Value2 -eq "Dbl") {
            $cell.Value = "SGL"
        }
    }
}

# Restore the author's last active selection on the sheet.
$ws.Range("C35").Select()
